# Update bases das ligas (10-06-2024 07:08)
# The edit swaps the full data payload (columns B:AD) between specific
# row pairs in the "Portugal Segunda Liga" sheet, while leaving column A
# (the sequential index) untouched on each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $range1 = $ws.Range("B$r1:AD$r1")
    $range2 = $ws.Range("B$r2:AD$r2")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}

Swap-Rows 77 78
Swap-Rows 132 133
Swap-Rows 140 141
Swap-Rows 230 231
